# Daily attendance processing - 2025-11-21 05:47:40
# Reorders the "Recorded By" (column G) list for each attendance row:
# whenever the comma-separated list of recorders begins with "System"
# (or "system") or "backup@backdoor.com", that leading entry is moved
# to the end of the list (a left-rotation by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Length -lt 2) { continue }

    $first = $parts[0].Trim()
    if ($first -eq "System" -or $first -eq "system" -or $first -eq "backup@backdoor.com") {
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + ,$first
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value2 = $newVal
    }
}
